# Apply scheduled-runner market data refresh to the Halicarnassus Profits workbook.
# Each sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) holds a Leve-profit table;
# columns H:N are live market-derived figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# that get refreshed with each run. No formulas are involved - every cell is a literal value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value2 = 39999.668
$ws.Range("J3").Value2 = 39999.668
$ws.Range("L3").Value2 = 39999.668
$ws.Range("N3").Value2 = -40227.668
$ws.Range("H12").Value2 = 145.4
$ws.Range("I12").Value2 = 80
$ws.Range("J12").Value2 = 161.75
$ws.Range("K12").Value2 = 80
$ws.Range("L12").Value2 = 161.75
$ws.Range("M12").Value2 = 90
$ws.Range("N12").Value2 = -501.75
$ws.Range("H18").Value2 = 769
$ws.Range("I18").Value2 = 769
$ws.Range("J18").Value2 = 0
$ws.Range("K18").Value2 = 769
$ws.Range("L18").Value2 = 0
$ws.Range("M18").Value2 = -485
$ws.Range("N18").ClearContents()
$ws.Range("H40").Value2 = 5264.95
$ws.Range("J40").Value2 = 6515.3335
$ws.Range("L40").Value2 = 6515.3335
$ws.Range("N40").Value2 = -6865.3335
$ws.Range("H55").Value2 = 84.59999999999999
$ws.Range("J55").Value2 = 78.333336
$ws.Range("L55").Value2 = 78.333336
$ws.Range("N55").Value2 = -506.333336
$ws.Range("H70").Value2 = 3552.6365
$ws.Range("I70").Value2 = 899
$ws.Range("J70").Value2 = 4142.3335
$ws.Range("K70").Value2 = 2697
$ws.Range("L70").Value2 = 12427.0005
$ws.Range("M70").Value2 = -2427
$ws.Range("N70").Value2 = -12967.0005
$ws.Range("H73").Value2 = 3552.6365
$ws.Range("I73").Value2 = 899
$ws.Range("J73").Value2 = 4142.3335
$ws.Range("K73").Value2 = 2697
$ws.Range("L73").Value2 = 12427.0005
$ws.Range("M73").Value2 = -1761
$ws.Range("N73").Value2 = -14299.0005
$ws.Range("H94").Value2 = 6433.643
$ws.Range("I94").Value2 = 6433.643
$ws.Range("K94").Value2 = 6433.643
$ws.Range("M94").Value2 = -5982.643
$ws.Range("H102").Value2 = 39999.668
$ws.Range("J102").Value2 = 39999.668
$ws.Range("L102").Value2 = 39999.668
$ws.Range("N102").Value2 = -46489.668
$ws.Range("H112").Value2 = 5001500
$ws.Range("J112").Value2 = 3000
$ws.Range("L112").Value2 = 9000
$ws.Range("N112").Value2 = -11216
$ws.Range("H137").Value2 = 3380.6667
$ws.Range("I137").Value2 = 2250
$ws.Range("J137").Value2 = 3815.5386
$ws.Range("K137").Value2 = 6750
$ws.Range("L137").Value2 = 11446.6158
$ws.Range("M137").Value2 = -4200
$ws.Range("N137").Value2 = -16546.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 3541.5557
$ws.Range("I45").Value2 = 2479
$ws.Range("K45").Value2 = 2479
$ws.Range("M45").Value2 = -2102
$ws.Range("H55").Value2 = 53999.5
$ws.Range("J55").Value2 = 99999
$ws.Range("L55").Value2 = 99999
$ws.Range("N55").Value2 = -100629
$ws.Range("H61").Value2 = 3933.5
$ws.Range("I61").Value2 = 3933.5
$ws.Range("J61").Value2 = 0
$ws.Range("K61").Value2 = 3933.5
$ws.Range("L61").Value2 = 0
$ws.Range("M61").Value2 = -3721.5
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value2 = 5609.6924
$ws.Range("I74").Value2 = 4812.636
$ws.Range("J74").Value2 = 9993.5
$ws.Range("K74").Value2 = 4812.636
$ws.Range("L74").Value2 = 9993.5
$ws.Range("M74").Value2 = -3938.636
$ws.Range("N74").Value2 = -11741.5
$ws.Range("H77").Value2 = 5609.6924
$ws.Range("I77").Value2 = 4812.636
$ws.Range("J77").Value2 = 9993.5
$ws.Range("K77").Value2 = 24063.18
$ws.Range("L77").Value2 = 49967.5
$ws.Range("M77").Value2 = -19695.18
$ws.Range("N77").Value2 = -58703.5
$ws.Range("H132").Value2 = 2089.125
$ws.Range("I132").Value2 = 1605.5
$ws.Range("K132").Value2 = 4816.5
$ws.Range("M132").Value2 = -2286.5
$ws.Range("H136").Value2 = 3933.5
$ws.Range("I136").Value2 = 3933.5
$ws.Range("J136").Value2 = 0
$ws.Range("K136").Value2 = 11800.5
$ws.Range("L136").Value2 = 0
$ws.Range("M136").Value2 = -9250.5
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 4611.6924
$ws.Range("J31").Value2 = 7224.4
$ws.Range("L31").Value2 = 7224.4
$ws.Range("N31").Value2 = -7814.4
$ws.Range("H34").Value2 = 4611.6924
$ws.Range("J34").Value2 = 7224.4
$ws.Range("L34").Value2 = 7224.4
$ws.Range("N34").Value2 = -7628.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value2 = 3383.875
$ws.Range("I139").Value2 = 2390.1667
$ws.Range("K139").Value2 = 7170.500100000001
$ws.Range("M139").Value2 = -2030.500100000001
$ws.Range("H141").Value2 = 1369.75
$ws.Range("I141").Value2 = 1369.75
$ws.Range("K141").Value2 = 4109.25
$ws.Range("M141").Value2 = 1070.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 1426.4286
$ws.Range("J22").Value2 = 1997.5
$ws.Range("L22").Value2 = 1997.5
$ws.Range("N22").Value2 = -2587.5
$ws.Range("H27").Value2 = 1426.4286
$ws.Range("J27").Value2 = 1997.5
$ws.Range("L27").Value2 = 1997.5
$ws.Range("N27").Value2 = -2211.5
$ws.Range("H68").Value2 = 9125
$ws.Range("I68").Value2 = 2800
$ws.Range("K68").Value2 = 2800
$ws.Range("M68").Value2 = -2051
$ws.Range("H71").Value2 = 9125
$ws.Range("I71").Value2 = 2800
$ws.Range("K71").Value2 = 14000
$ws.Range("M71").Value2 = -10256
$ws.Range("H122").Value2 = 3386.6667
$ws.Range("I122").Value2 = 3386.6667
$ws.Range("K122").Value2 = 10160.0001
$ws.Range("M122").Value2 = -7710.000100000001
$ws.Range("H132").Value2 = 5345.091
$ws.Range("I132").Value2 = 4310.778
$ws.Range("K132").Value2 = 12932.334
$ws.Range("M132").Value2 = -10402.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 7552.636
$ws.Range("I62").Value2 = 5246.75
$ws.Range("J62").Value2 = 8870.286
$ws.Range("K62").Value2 = 5246.75
$ws.Range("L62").Value2 = 8870.286
$ws.Range("M62").Value2 = -4622.75
$ws.Range("N62").Value2 = -10118.286
$ws.Range("H65").Value2 = 7552.636
$ws.Range("I65").Value2 = 5246.75
$ws.Range("J65").Value2 = 8870.286
$ws.Range("K65").Value2 = 26233.75
$ws.Range("L65").Value2 = 44351.43
$ws.Range("M65").Value2 = -23113.75
$ws.Range("N65").Value2 = -50591.43
$ws.Range("H96").Value2 = 0
$ws.Range("I96").Value2 = 0
$ws.Range("J96").Value2 = 0
$ws.Range("K96").Value2 = 0
$ws.Range("L96").Value2 = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()
$ws.Range("H126").Value2 = 2962.724
$ws.Range("I126").Value2 = 1416.0526
$ws.Range("J126").Value2 = 5901.4
$ws.Range("K126").Value2 = 4248.1578
$ws.Range("L126").Value2 = 17704.2
$ws.Range("M126").Value2 = -1778.1578
$ws.Range("N126").Value2 = -22644.2
$ws.Range("H132").Value2 = 2097.8823
$ws.Range("I132").Value2 = 2043.4667
$ws.Range("K132").Value2 = 6130.4001
$ws.Range("M132").Value2 = -3600.4001
